# Reporte de inventario - ajustes de ancho de columna e insercion de columna ID + nueva fila
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the generated-at timestamp in A1 (merged A1:F1)
$ws.Range("A1").Value2 = "Reporte generado el 19/05/2025 a las 11:56"

# 2. Insert a new column before column B ("Nombre") to host the new "ID" column.
#    This shifts Nombre/Cantidad/Precio/Unidad/Categoria one column to the right,
#    duplicating formatting from column A into the new column B automatically.
$ws.Columns.Item(2).Insert()

# Column insert auto-extends the merged title range to A1:G1; restore it to A1:F1
$ws.Range("A1:G1").UnMerge()
$ws.Range("A1:F1").Merge()

# 3. Header row (row 2): set the new "ID" header text (style already matches the
#    rest of the green header band because it was copied from column A on insert)
$ws.Range("B2").Value2 = "ID"

# The "N°." header cell (A2) switches to a blue fill, everything else stays the same
$ws.Range("A2").Interior.Color = 13792793

# 4. Fill in the "ID" values for the existing data rows (style already correct,
#    copied automatically from column A when the column was inserted)
$ws.Range("B3").Value2 = 62
$ws.Range("B4").Value2 = 63

# 5. Add the new inventory row (row 5), cloning the formatting of row 4 first so
#    no extra/duplicate style entries are generated, then filling in the values
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A5").Value2 = 3
$ws.Range("B5").Value2 = 64
$ws.Range("C5").Value2 = "Oeoeoeoeoeoeoeoeoeoe"
$ws.Range("D5").Value2 = 12
$ws.Range("E5").Value2 = 12
$ws.Range("F5").Value2 = "par"
$ws.Range("G5").Value2 = "Sala"

# 6. Column widths (character units). Excel stores widths with a small built-in
#    padding offset (~0.8333 chars for the default Calibri 11 font), so subtract
#    that offset to land exactly on the desired stored width.
$pad = 0.8333333333333
$ws.Columns.Item(1).ColumnWidth = 10 - $pad
$ws.Columns.Item(2).ColumnWidth = 10 - $pad
$ws.Columns.Item(3).ColumnWidth = 35 - $pad
$ws.Columns.Item(4).ColumnWidth = 26 - $pad
$ws.Columns.Item(5).ColumnWidth = 26 - $pad
$ws.Columns.Item(6).ColumnWidth = 26 - $pad
$ws.Columns.Item(7).ColumnWidth = 25 - $pad

Write-Host "edit complete"
